# Automatic update of files.
# Bump the "Förändrad" (Changed) date column C from 45181 to 45182
# for every data row (rows 2-33) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45181) {
        $cell.Value2 = 45182
    }
}
